$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 219.72728
$ws.Range("I5").Value = 219.72728
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 219.72728
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -104.72728
$ws.Range("N5").ClearContents()
$ws.Range("H12").Value = 100
$ws.Range("I12").Value = 100
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 100
$ws.Range("L12").Value = 100
$ws.Range("M12").Value = 70
$ws.Range("N12").Value = -440
$ws.Range("H70").Value = 885.2857
$ws.Range("I70").Value = 832.0833
$ws.Range("J70").Value = 925.1875
$ws.Range("K70").Value = 2496.2499
$ws.Range("L70").Value = 2775.5625
$ws.Range("M70").Value = -2226.2499
$ws.Range("N70").Value = -3315.5625
$ws.Range("H73").Value = 885.2857
$ws.Range("I73").Value = 832.0833
$ws.Range("J73").Value = 925.1875
$ws.Range("K73").Value = 2496.2499
$ws.Range("L73").Value = 2775.5625
$ws.Range("M73").Value = -1560.2499
$ws.Range("N73").Value = -4647.5625
$ws.Range("H116").Value = 4365.919
$ws.Range("I116").Value = 4531
$ws.Range("J116").Value = 4209.5264
$ws.Range("K116").Value = 4531
$ws.Range("L116").Value = 4209.5264
$ws.Range("M116").Value = -1089
$ws.Range("N116").Value = -11093.5264
$ws.Range("H129").Value = 704.1429000000001
$ws.Range("I129").Value = 416.34784
$ws.Range("J129").Value = 2028
$ws.Range("K129").Value = 1249.04352
$ws.Range("L129").Value = 6084
$ws.Range("M129").Value = 3750.95648
$ws.Range("N129").Value = -16084
$ws.Range("H134").Value = 45249.5
$ws.Range("J134").Value = 45249.5
$ws.Range("L134").Value = 45249.5
$ws.Range("N134").Value = -55389.5
$ws.Range("H137").Value = 2273704.2
$ws.Range("I137").Value = 981331.0600000001
$ws.Range("J137").Value = 6667773.5
$ws.Range("K137").Value = 2943993.18
$ws.Range("L137").Value = 20003320.5
$ws.Range("M137").Value = -2941443.18
$ws.Range("N137").Value = -20008420.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20958.148
$ws.Range("I32").Value = 5193.3237
$ws.Range("K32").Value = 5193.3237
$ws.Range("M32").Value = -4906.3237
$ws.Range("H88").Value = 22129.1
$ws.Range("I88").Value = 1535.2
$ws.Range("J88").Value = 42723
$ws.Range("K88").Value = 1535.2
$ws.Range("L88").Value = 42723
$ws.Range("M88").Value = -1129.2
$ws.Range("N88").Value = -43535
$ws.Range("H91").Value = 22129.1
$ws.Range("I91").Value = 1535.2
$ws.Range("J91").Value = 42723
$ws.Range("K91").Value = 1535.2
$ws.Range("L91").Value = 42723
$ws.Range("M91").Value = -131.2
$ws.Range("N91").Value = -45531
$ws.Range("H132").Value = 121310.88
$ws.Range("I132").Value = 144333.52
$ws.Range("J132").Value = 6197.7144
$ws.Range("K132").Value = 433000.5599999999
$ws.Range("L132").Value = 18593.1432
$ws.Range("M132").Value = -430470.5599999999
$ws.Range("N132").Value = -23653.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 929.35
$ws.Range("I64").Value = 958.8570999999999
$ws.Range("J64").Value = 913.46155
$ws.Range("K64").Value = 958.8570999999999
$ws.Range("L64").Value = 913.46155
$ws.Range("M64").Value = -733.8570999999999
$ws.Range("N64").Value = -1363.46155
$ws.Range("H67").Value = 929.35
$ws.Range("I67").Value = 958.8570999999999
$ws.Range("J67").Value = 913.46155
$ws.Range("K67").Value = 958.8570999999999
$ws.Range("L67").Value = 913.46155
$ws.Range("M67").Value = -178.8570999999999
$ws.Range("N67").Value = -2473.46155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3111.111
$ws.Range("I62").Value = 2857.1428
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 2857.1428
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -2233.1428
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 3111.111
$ws.Range("I65").Value = 2857.1428
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 14285.714
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -11165.714
$ws.Range("N65").Value = -26240

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1971449.9
$ws.Range("I2").Value = 76946.92
$ws.Range("J2").Value = 14285719
$ws.Range("K2").Value = 461681.52
$ws.Range("L2").Value = 85714314
$ws.Range("M2").Value = -461568.52
$ws.Range("N2").Value = -85714540
$ws.Range("H82").Value = 6856.5386
$ws.Range("I82").Value = 750.4
$ws.Range("J82").Value = 8310.380999999999
$ws.Range("K82").Value = 2251.2
$ws.Range("L82").Value = 24931.143
$ws.Range("M82").Value = -1845.2
$ws.Range("N82").Value = -25743.143
$ws.Range("H85").Value = 6856.5386
$ws.Range("I85").Value = 750.4
$ws.Range("J85").Value = 8310.380999999999
$ws.Range("K85").Value = 2251.2
$ws.Range("L85").Value = 24931.143
$ws.Range("M85").Value = -847.1999999999998
$ws.Range("N85").Value = -27739.143
$ws.Range("H88").Value = 3887.5
$ws.Range("J88").Value = 3887.5
$ws.Range("L88").Value = 11662.5
$ws.Range("N88").Value = -12518.5
$ws.Range("H91").Value = 3887.5
$ws.Range("J91").Value = 3887.5
$ws.Range("L91").Value = 11662.5
$ws.Range("N91").Value = -14626.5
$ws.Range("H94").Value = 5189
$ws.Range("J94").Value = 4935.7896
$ws.Range("L94").Value = 14807.3688
$ws.Range("N94").Value = -16159.3688
$ws.Range("H131").Value = 972.9677
$ws.Range("J131").Value = 1009.02325
$ws.Range("L131").Value = 3027.06975
$ws.Range("N131").Value = -13107.06975

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 21373.242
$ws.Range("I70").Value = 25035.02
$ws.Range("J70").Value = 4395.909
$ws.Range("K70").Value = 25035.02
$ws.Range("L70").Value = 4395.909
$ws.Range("M70").Value = -24765.02
$ws.Range("N70").Value = -4935.909
$ws.Range("H73").Value = 21373.242
$ws.Range("I73").Value = 25035.02
$ws.Range("J73").Value = 4395.909
$ws.Range("K73").Value = 25035.02
$ws.Range("L73").Value = 4395.909
$ws.Range("M73").Value = -24099.02
$ws.Range("N73").Value = -6267.909
$ws.Range("H113").Value = 1603
$ws.Range("I113").Value = 1392.5
$ws.Range("J113").Value = 1843.5714
$ws.Range("K113").Value = 1392.5
$ws.Range("L113").Value = 1843.5714
$ws.Range("M113").Value = 777.5
$ws.Range("N113").Value = -6183.5714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2032.8182
$ws.Range("I93").Value = 1720.6
$ws.Range("J93").Value = 2293
$ws.Range("K93").Value = 1720.6
$ws.Range("L93").Value = 2293
$ws.Range("M93").Value = -472.5999999999999
$ws.Range("N93").Value = -4789
$ws.Range("H132").Value = 1793.0754
$ws.Range("I132").Value = 1671.5209
$ws.Range("K132").Value = 5014.5627
$ws.Range("M132").Value = -2484.5627

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 476882.84
$ws.Range("I100").Value = 468.7
$ws.Range("J100").Value = 909986.6
$ws.Range("K100").Value = 937.4
$ws.Range("L100").Value = 1819973.2
$ws.Range("M100").Value = -396.4
$ws.Range("N100").Value = -1821055.2
$ws.Range("H136").Value = 1480.9344
$ws.Range("I136").Value = 1525.5625
$ws.Range("J136").Value = 1316.1538
$ws.Range("K136").Value = 4576.6875
$ws.Range("L136").Value = 3948.4614
$ws.Range("M136").Value = -2026.6875
$ws.Range("N136").Value = -9048.4614
